# Update to the latest Business count data.
# Rows 8 and 9 hold the "ONS UK Business Counts" source rows
# (Enterprises by employment size band / Enterprises by employment industry).
# Column C = "Latest period (release date)", column D = "Next period (release date)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value = "Mar 2025 (14/10/24)"
$ws.Range("D8").Value = "Mar 2026 (Autumn 26)"
$ws.Range("C9").Value = "Mar 2025 (14/10/24)"
$ws.Range("D9").Value = "Mar 2026 (Autumn 26)"

# Restore the selection state that was recorded when the file was last saved.
$ws.Range("E11").Select()
